$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Company name header (A2) changed from "Bottlejac Trading Pty Ltd"
#    to "CP & SA Webster"
# ------------------------------------------------------------------
$ws.Range("A2").Value = "CP & SA Webster"

# ------------------------------------------------------------------
# 2. Updated "Cattle" column (C) figures
# ------------------------------------------------------------------
$ws.Range("C6").Value  = 1308
$ws.Range("C8").Value  = 23
$ws.Range("C11").Value = 1530
$ws.Range("C14").Value = 853

# ------------------------------------------------------------------
# 3. New "Horses" category added in column M, mirroring the layout
#    already used for the other livestock categories (Cattle/Horse/
#    Goats/Pigs/Rams in columns C/E/G/I/K). Copy formatting from the
#    equivalent Cattle (C) cell for every row, then set the new value.
# ------------------------------------------------------------------
$ws.Range("C4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = "Horses"

$ws.Range("C5").Copy()
$ws.Range("M5").PasteSpecial(-4122)
$ws.Range("M5").Value = "No."

$ws.Range("C6").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 21

$ws.Range("C8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M8").Value = 1

$ws.Range("C11").Copy()
$ws.Range("M11").PasteSpecial(-4122)
$ws.Range("M11").Value = 0

$ws.Range("C14").Copy()
$ws.Range("M14").PasteSpecial(-4122)
$ws.Range("M14").Value = 0

$ws.Range("C16").Copy()
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("M16").Value = 0

$ws.Range("C18").Copy()
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("M18").Value = 0

$ws.Range("C20").Copy()
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("M20").Formula = "=M6+M8+M11-M14-M16-M18"

# Give column M (Horses) the same custom width as the other data
# columns (as close as this host's width-rounding allows to 11.125).
$ws.Columns.Item(13).ColumnWidth = 10.41

# ------------------------------------------------------------------
# 4. Cosmetic: move the active cell selection, matching where the
#    author left the cursor after editing.
# ------------------------------------------------------------------
[void]$ws.Range("S17").Select()
